$wb = $excel.ActiveWorkbook

$sheetNames = @("展览", "全部类型")

foreach ($name in $sheetNames) {
    $ws = $wb.Worksheets.Item($name)
    $ws.Range("F2").Value = 4725
    $ws.Range("F3").Value = 140
    $ws.Range("F4").Value = 824
}
